$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet title (reflects new "through" date)
$ws.Name = "Through 2022-06-02"

# Update the label in A7 ("June (through 06-01)" -> "June (through 06-02)")
$ws.Range("A7").Value = "June (through 06-02)"

# Update June row (row 7) values
$ws.Range("C7").Value = 4
$ws.Range("E7").Value = 8
$ws.Range("F7").Value = 3
$ws.Range("G7").Value = 20
$ws.Range("H7").Value = 4
$ws.Range("I7").Value = 6

# Update Total row (row 8) values
$ws.Range("C8").Value = 213
$ws.Range("E8").Value = 303
$ws.Range("F8").Value = 207
$ws.Range("G8").Value = 378
$ws.Range("H8").Value = 635
$ws.Range("I8").Value = 670
